$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This paragraph's text originally reads (single run per segment):
#   "The solder mask does not do well when exposed quickly in my custom uiv rig. ..."
# The edit:
#   1) splits "...when exposed..." into two runs ("...when" | " exposed...custom ")
#      with no textual change, and
#   2) fixes the typo "uiv" -> "uv", split into two runs ("u" | "v").
# Both new runs keep the original (bold, blue, 36pt) character formatting.
# ---------------------------------------------------------------------------

# --- Step 1: fix the "uiv" -> "uv" typo -------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("uiv")
$rTypo = $d.Range($idx, $idx + 3)
$rTypo.Text = "uv"

# --- Step 2: split "The solder mask ... when" | " exposed ... custom " ------------------
$full = $d.Content.Text
$prefix = "The solder mask does not do well when"
$splitPos = $full.IndexOf($prefix) + $prefix.Length
$suffixLen = (" exposed quickly in my custom ").Length
$rSplit1 = $d.Range($splitPos, $splitPos + $suffixLen)
# Toggling Bold off then back on forces the run to be split at this boundary
# while leaving the (identical) resulting formatting untouched.
$rSplit1.Bold = 0
$rSplit1.Bold = 1

# --- Step 3: split "u" | "v" --------------------------------------------------------------
$full = $d.Content.Text
$idx2 = $full.IndexOf("uv rig")
$rSplit2 = $d.Range($idx2 + 1, $idx2 + 2)
$rSplit2.Bold = 0
$rSplit2.Bold = 1
